$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.030832953909584
$ws.Range("D2").Value = 1.04154846395327
$ws.Range("E2").Value = 1.030489117611669
$ws.Range("F2").Value = 1.052339690071501
$ws.Range("I2").Value = 1.038538951906923
$ws.Range("J2").Value = 1.035972018432981
$ws.Range("K2").Value = 1.044327666703657
$ws.Range("L2").Value = 1.033299981444586
$ws.Range("M2").Value = 1.055088694902679
$ws.Range("N2").Value = 1.01598494316569
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.03160493092163
$ws.Range("D3").Value = 1.042169336270755
$ws.Range("E3").Value = 1.031140494892122
$ws.Range("F3").Value = 1.053127480209819
$ws.Range("I3").Value = 1.038712502810314
$ws.Range("J3").Value = 1.036386520944202
$ws.Range("K3").Value = 1.044759687535406
$ws.Range("L3").Value = 1.033760146016331
$ws.Range("M3").Value = 1.055689375381174
$ws.Range("N3").Value = 1.016123405705688
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.032105203430294
$ws.Range("D4").Value = 1.042571755423777
$ws.Range("E4").Value = 1.031563012118707
$ws.Range("F4").Value = 1.053638302146798
$ws.Range("I4").Value = 1.038824021425752
$ws.Range("J4").Value = 1.036654797599265
$ws.Range("K4").Value = 1.045039213645047
$ws.Range("L4").Value = 1.034058249533407
$ws.Range("M4").Value = 1.056078464758798
$ws.Range("N4").Value = 1.01621299114758
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.032315696053148
$ws.Range("D5").Value = 1.042741091645451
$ws.Range("E5").Value = 1.031740883575571
$ws.Range("F5").Value = 1.053853305417767
$ws.Range("I5").Value = 1.038870716226649
$ws.Range("J5").Value = 1.036767595527379
$ws.Range("K5").Value = 1.045156720122057
$ws.Range("L5").Value = 1.034183653680741
$ws.Range("M5").Value = 1.056242134006925
$ws.Range("N5").Value = 1.016250650210453
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.032351049076695
$ws.Range("D6").Value = 1.042769533221081
$ws.Range("E6").Value = 1.031770763297603
$ws.Range("F6").Value = 1.053889420227764
$ws.Range("I6").Value = 1.038878545463829
$ws.Range("J6").Value = 1.036786535602795
$ws.Range("K6").Value = 1.045176449565029
$ws.Range("L6").Value = 1.034204714321173
$ws.Range("M6").Value = 1.056269620374526
$ws.Range("N6").Value = 1.01625697316296
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.032108015345763
$ws.Range("D7").Value = 1.042574017479549
$ws.Range("E7").Value = 1.031565387884792
$ws.Range("F7").Value = 1.053641174036813
$ws.Range("I7").Value = 1.038824646102189
$ws.Range("J7").Value = 1.036656304755705
$ws.Range("K7").Value = 1.045040783798916
$ws.Range("L7").Value = 1.034059924871648
$ws.Range("M7").Value = 1.056080651340092
$ws.Range("N7").Value = 1.016213494360868
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.031093690347325
$ws.Range("D8").Value = 1.041758150077691
$ws.Range("E8").Value = 1.030709038620623
$ws.Range("F8").Value = 1.052605704762043
$ws.Range("I8").Value = 1.038597765363863
$ws.Range("J8").Value = 1.036112087068614
$ws.Range("K8").Value = 1.044473673489148
$ws.Range("L8").Value = 1.033455423636187
$ws.Range("M8").Value = 1.055291611909925
$ws.Range("N8").Value = 1.016031738819176
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.029312162420187
$ws.Range("D9").Value = 1.040325730714888
$ws.Range("E9").Value = 1.029208037789819
$ws.Range("F9").Value = 1.050789360581854
$ws.Range("I9").Value = 1.038192030413034
$ws.Range("J9").Value = 1.035153669020431
$ws.Range("K9").Value = 1.043474261674146
$ws.Range("L9").Value = 1.032392932181632
$ws.Range("M9").Value = 1.053904436500924
$ws.Range("N9").Value = 1.015711414147361
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.028128519707042
$ws.Range("D10").Value = 1.039374427448117
$ws.Range("E10").Value = 1.028212865393521
$ws.Range("F10").Value = 1.04958416665052
$ws.Range("I10").Value = 1.037917592704554
$ws.Range("J10").Value = 1.034515182027132
$ws.Range("K10").Value = 1.042808012756944
$ws.Range("L10").Value = 1.031686519166907
$ws.Range("M10").Value = 1.052981922494564
$ws.Range("N10").Value = 1.015497860858315
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.027616972441739
$ws.Range("D11").Value = 1.038963391682741
$ws.Range("E11").Value = 1.027783273283039
$ws.Range("F11").Value = 1.049063684833333
$ws.Range("I11").Value = 1.03779783183638
$ws.Range("J11").Value = 1.034238835748236
$ws.Range("K11").Value = 1.04251954410532
$ws.Range("L11").Value = 1.031381107302203
$ws.Range("M11").Value = 1.05258302527904
$ws.Range("N11").Value = 1.01540539521777
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.027427109887987
$ws.Range("D12").Value = 1.038810849671759
$ws.Range("E12").Value = 1.027623904641162
$ws.Range("F12").Value = 1.048870563589304
$ws.Range("I12").Value = 1.037753208769721
$ws.Range("J12").Value = 1.03413620825768
$ws.Range("K12").Value = 1.04241239873056
$ws.Range("L12").Value = 1.031267735812091
$ws.Range("M12").Value = 1.052434942728683
$ws.Range("N12").Value = 1.015371050510843
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.02746782928165
$ws.Range("D13").Value = 1.038843564310394
$ws.Range("E13").Value = 1.027658080627266
$ws.Range("F13").Value = 1.048911979262937
$ws.Range("I13").Value = 1.037762786825672
$ws.Range("J13").Value = 1.034158221283785
$ws.Range("K13").Value = 1.042435381539321
$ws.Range("L13").Value = 1.031292051095675
$ws.Range("M13").Value = 1.052466703032886
$ws.Range("N13").Value = 1.015378417507104
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.02760127529538
$ws.Range("D14").Value = 1.03895077974733
$ws.Range("E14").Value = 1.027770095709217
$ws.Range("F14").Value = 1.049047717096826
$ws.Range("I14").Value = 1.037794146101454
$ws.Range("J14").Value = 1.034230352115562
$ws.Range("K14").Value = 1.04251068733017
$ws.Range("L14").Value = 1.031371734500386
$ws.Range("M14").Value = 1.052570782971363
$ws.Range("N14").Value = 1.015402556247273
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.027683515575282
$ws.Range("D15").Value = 1.039016856662816
$ws.Range("E15").Value = 1.027839138593696
$ws.Range("F15").Value = 1.049131377402729
$ws.Range("I15").Value = 1.037813449251359
$ws.Range("J15").Value = 1.034274796969239
$ws.Range("K15").Value = 1.042557086384255
$ws.Range("L15").Value = 1.031420839671908
$ws.Range("M15").Value = 1.052634921460633
$ws.Range("N15").Value = 1.015417429088865
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.028162490184197
$ws.Range("D16").Value = 1.039401725354505
$ws.Range("E16").Value = 1.028241404080678
$ws.Range("F16").Value = 1.049618738481409
$ws.Range("I16").Value = 1.037925521365476
$ws.Range("J16").Value = 1.03453352492398
$ws.Range("K16").Value = 1.042827158052248
$ws.Range("L16").Value = 1.031706798370745
$ws.Range("M16").Value = 1.053008407894739
$ws.Range("N16").Value = 1.015503997632247
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.02846320148616
$ws.Range("D17").Value = 1.039643381649915
$ws.Range("E17").Value = 1.028494090551367
$ws.Range("F17").Value = 1.04992481697765
$ws.Range("I17").Value = 1.037995573482124
$ws.Range("J17").Value = 1.034695851981004
$ws.Range("K17").Value = 1.042996573691336
$ws.Range("L17").Value = 1.031886299513636
$ws.Range("M17").Value = 1.053242836642082
$ws.Range("N17").Value = 1.015558301305822
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.028638695516312
$ws.Range("D18").Value = 1.039784420914945
$ws.Range("E18").Value = 1.028641605913258
$ws.Range("F18").Value = 1.050103479911044
$ws.Range("I18").Value = 1.038036344149508
$ws.Range("J18").Value = 1.034790546363701
$ws.Range("K18").Value = 1.043095392944237
$ws.Range("L18").Value = 1.031991044704371
$ws.Range("M18").Value = 1.053379628622207
$ws.Range("N18").Value = 1.01558997610563
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.02869855036614
$ws.Range("D19").Value = 1.039832526055053
$ws.Range("E19").Value = 1.028691926423018
$ws.Range("F19").Value = 1.05016442176333
$ws.Range("I19").Value = 1.038050230680336
$ws.Range("J19").Value = 1.034822836653465
$ws.Range("K19").Value = 1.04312908804566
$ws.Range("L19").Value = 1.032026767729915
$ws.Range("M19").Value = 1.053426280196295
$ws.Range("N19").Value = 1.01560077644413
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.028430928250044
$ws.Range("D20").Value = 1.039617445385332
$ws.Range("E20").Value = 1.028466966469387
$ws.Range("F20").Value = 1.049891963927565
$ws.Range("I20").Value = 1.037988066809464
$ws.Range("J20").Value = 1.034678434598092
$ws.Range("K20").Value = 1.042978396781357
$ws.Range("L20").Value = 1.031867036052456
$ws.Range("M20").Value = 1.053217679092395
$ws.Range("N20").Value = 1.015552474992092
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.027561974644841
$ws.Range("D21").Value = 1.03891920373114
$ws.Range("E21").Value = 1.027737104494469
$ws.Range("F21").Value = 1.04900773990356
$ws.Range("I21").Value = 1.03778491539449
$ws.Range("J21").Value = 1.034209110825226
$ws.Range("K21").Value = 1.042488511509257
$ws.Range("L21").Value = 1.031348267721955
$ws.Range("M21").Value = 1.052540131644467
$ws.Range("N21").Value = 1.015395447954233
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.027016491357679
$ws.Range("D22").Value = 1.038480973412604
$ws.Range("E22").Value = 1.027279375219727
$ws.Range("F22").Value = 1.048453002793626
$ws.Range("I22").Value = 1.037656384840167
$ws.Range("J22").Value = 1.033914143702942
$ws.Range("K22").Value = 1.042180529279563
$ws.Range("L22").Value = 1.031022515235892
$ws.Range("M22").Value = 1.05211462707166
$ws.Range("N22").Value = 1.015296725716248
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.027305579972314
$ws.Range("D23").Value = 1.038713212798998
$ws.Range("E23").Value = 1.027521915154766
$ws.Range("F23").Value = 1.04874696401086
$ws.Range("I23").Value = 1.037724596994563
$ws.Range("J23").Value = 1.034070499962568
$ws.Range("K23").Value = 1.042343793326518
$ws.Range("L23").Value = 1.03119516268764
$ws.Range("M23").Value = 1.052340147413465
$ws.Range("N23").Value = 1.01534905942347
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.028445510845411
$ws.Range("D24").Value = 1.039629164603795
$ws.Range("E24").Value = 1.028479222280455
$ws.Range("F24").Value = 1.049906808397585
$ws.Range("I24").Value = 1.037991459028511
$ws.Range("J24").Value = 1.034686304727566
$ws.Range("K24").Value = 1.042986610138653
$ws.Range("L24").Value = 1.031875740241889
$ws.Range("M24").Value = 1.053229046540665
$ws.Range("N24").Value = 1.015555107651351
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.029772025954283
$ws.Range("D25").Value = 1.040695411964735
$ws.Range("E25").Value = 1.029595122667179
$ws.Range("F25").Value = 1.051257933688918
$ws.Range("I25").Value = 1.03829762191572
$ws.Range("J25").Value = 1.035401368109451
$ws.Range("K25").Value = 1.043732635420937
$ws.Range("L25").Value = 1.032667280654844
$ws.Range("M25").Value = 1.054262662572873
$ws.Range("N25").Value = 1.015794228484633
